# Regenerate save_data to use K instead of Strike# (recalculated K values)
# Updates column G (K) values for rows 2-10 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 2
    4  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 3
    10 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
